$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Insert a new row at position 11 (after "Contact", before "Description"),
# shifting everything below it down by one row.
$ws.Rows.Item(11).Insert()

# Populate the newly inserted row with the "Jurisdiction" property (empty value).
$ws.Cells.Item(11, 1).Value = "Jurisdiction"
$ws.Cells.Item(11, 2).Value = ""

# Update the publication Date value (row 8, now still row 8 since insert was below it).
$ws.Cells.Item(8, 2).Value = "2024-07-01T07:50:29+00:00"
